$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 343 (existing rows 343:371 shift down to 344:372)
$ws.Rows.Item(343).Insert()

# Populate the newly inserted row 343 with the new weekly record
$ws.Range("A343").Value = 6
$ws.Range("B343").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C343").Value = "Metropolitana"
$ws.Range("D343").Value = 45223
$ws.Range("E343").Value = 13
$ws.Range("F343").Value = 100112029
$ws.Range("G343").Value = "Orégano"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 33
$ws.Range("K343").Value = 16000
$ws.Range("L343").Value = 16000
$ws.Range("M343").Value = 16000
$ws.Range("N343").Value = "`$/docena de atados"
$ws.Range("O343").Value = "Región Metropolitana"
$ws.Range("P343").Value = 5333
$ws.Range("Q343").Value = 3
$ws.Range("R343").Value = "Hortaliza"
